$d = $word.ActiveDocument

# 1. Merge "Cramer Fish " + "Sciences (CFS), " into a single run.
$d.Content.Find.Execute(
    "Cramer Fish Sciences (CFS), ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Cramer Fish Sciences (CFS), ",
    2) | Out-Null

# 2. Merge the Chinook Salmon / steelhead sentence (several runs, two of
#    them italicized for the scientific names) into a single, non-italic
#    run. Word's Find/Replace on a range spanning multiple runs collapses
#    them into one run that uses the formatting of the first run in the
#    matched range (non-italic "normaltextrun"), which is exactly what the
#    diff shows (the italics on the species names are lost).
$d.Content.Find.Execute(
    " Side Channel and Floodplain Restoration Project (Project) on the Yuba River, California. The Project is designed to restore and enhance ecosystem processes, with a primary focus on improving productive juvenile salmonid rearing habitat to increase natural production of fall and spring-run Chinook Salmon (Oncorhynchus tshawytscha) and steelhead (O. mykiss) in the Yuba River. The Project would enhance and/or create up to 1",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    " Side Channel and Floodplain Restoration Project (Project) on the Yuba River, California. The Project is designed to restore and enhance ecosystem processes, with a primary focus on improving productive juvenile salmonid rearing habitat to increase natural production of fall and spring-run Chinook Salmon (Oncorhynchus tshawytscha) and steelhead (O. mykiss) in the Yuba River. The Project would enhance and/or create up to 1",
    2) | Out-Null

# 3. Merge " miles of seasonal side channels..." + "in order to" + " allow
#    natural river..." into a single run, removing the grammar-check
#    proofErr markers around "in order to" in the process.
$d.Content.Find.Execute(
    " miles of seasonal side channels. The design approach focuses on removing unnatural constraints (such as a mid-river training wall and very coarse surface materials left from mining activities) in order to allow natural river and floodplain processes to function. Construction planning efforts include multi-year phasing to remove ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    " miles of seasonal side channels. The design approach focuses on removing unnatural constraints (such as a mid-river training wall and very coarse surface materials left from mining activities) in order to allow natural river and floodplain processes to function. Construction planning efforts include multi-year phasing to remove ",
    2) | Out-Null
